# Append a summary "AVG" row under the existing run results (rows 2-41
# hold runs 1-40 in column A / accuracy in column B), mirroring the author's
# manual addition of an average row at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("A42").Value = "AVG"
$ws.Range("B42").Formula = "=AVERAGE(B2:B41)"

# Leave the freshly-entered average cell selected, scrolled into view -
# matches the saved sheetView (activeCell/sqref = B42, topLeftCell = A23).
$ws.Range("B42").Select()
$excel.ActiveWindow.ScrollRow = 23
